# Rework preprocessing to use a hashmap to prevent duplicate values.
#
# Column K holds Java-ish code template strings such as
#   if (fieldList[3].equals("1")) foods.add("skim milk");
# Column L substitutes the literal "skim milk" with the real food name
# from column D via =SUBSTITUTE(K#,"skim milk",D#).
#
# The preprocessing snippet was reworked to populate a hashmap instead of
# a list, so every occurrence of the `foods.add("skim milk");` call needs
# to become `foods.put("skim milk", 1);`. Updating column K is sufficient;
# column L recalculates automatically because it is a formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldCall = 'foods.add("skim milk");'
$newCall = 'foods.put("skim milk", 1);'

$dims = $ws.UsedRange
$lastRow = $dims.Rows.Count + $dims.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 11)
    $k = $cell.Value2
    if ($k -ne $null -and $k -ne "" -and $k.Contains($oldCall)) {
        $cell.Value2 = $k.Replace($oldCall, $newCall)
    }
}
